$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1 "Vendors": update inventory quantities (decrement/populate Quantity
# column) and move the current selection.
# ---------------------------------------------------------------------------
$wsVendors = $wb.Worksheets.Item(1)

$wsVendors.Cells.Item(2, 4).Value = 100      # D2  120 -> 100
$wsVendors.Cells.Item(5, 4).Value = 100      # D5  (blank) -> 100
$wsVendors.Cells.Item(9, 4).Value = 100      # D9  (blank) -> 100

# ---------------------------------------------------------------------------
# Sheet 2 "Clients": lay out the full client list used by the Welcome screen.
# ---------------------------------------------------------------------------
$wsClients = $wb.Worksheets.Item(2)

$clientNames = @(
  "Rosa",
  "Cheddar",
  "GZA",
  "asdf",
  "Doe",
  "qwerty",
  "zxcv",
  "wert",
  "ty",
  "fgh",
  "Guy",
  "Scully",
  "JP",
  "James",
  "GUy",
  "Garrison",
  "Gha",
  "Fuh",
  "Jameson",
  "Gus",
  "Hue Jass",
  "Faygo",
  "asdf",
  "Roger"
)

$row = 2
foreach ($name in $clientNames) {
  $wsClients.Cells.Item($row, 1).Value = 1
  $wsClients.Cells.Item($row, 2).Value = $name
  $row = $row + 1
}

# Row 26 holds a plain number instead of a name.
$wsClients.Cells.Item(26, 1).Value = 1
$wsClients.Cells.Item(26, 2).Value = 2

$wsClients.Cells.Item(27, 1).Value = 1
$wsClients.Cells.Item(27, 2).Value = "Hue Jass"

$wsClients.Cells.Item(28, 1).Value = 1
$wsClients.Cells.Item(28, 2).Value = "Qubert"

# ---------------------------------------------------------------------------
# Sheet 3 "Expenses": replace the sample data with the real header row used
# to drive the Welcome screen table (ClientID / ClientName / ProductName /
# ProductQuantity), and widen the columns to fit.
# ---------------------------------------------------------------------------
$wsExpenses = $wb.Worksheets.Item(3)

$wsExpenses.Cells.Clear()

$wsExpenses.Range("A1:D1").Font.Bold = $true

$wsExpenses.Cells.Item(1, 1).Value = "ClientID"
$wsExpenses.Cells.Item(1, 2).Value = "ClientName"
$wsExpenses.Cells.Item(1, 3).Value = "ProductName"
$wsExpenses.Cells.Item(1, 4).Value = "ProductQuantity"

$wsExpenses.Columns.Item(1).ColumnWidth = 7.333333333333333
$wsExpenses.Columns.Item(2).ColumnWidth = 10.333333333333332
$wsExpenses.Columns.Item(3).ColumnWidth = 12.0
$wsExpenses.Columns.Item(4).ColumnWidth = 14.833333333333332

$wsExpenses.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Selections / active sheet - restore the selections shown in the final
# workbook, finishing on the Vendors tab so it stays the active sheet.
# ---------------------------------------------------------------------------
$wsExpenses.Activate()
$wsExpenses.Columns.Item(1).Select()
$excel.ActiveWindow.Zoom = 130

$wsClients.Activate()
$wsClients.Range("D8").Select()

$wsVendors.Activate()
$wsVendors.Range("F8").Select()
